# Insert a new weekly price record as row 151 in the Mango sheet,
# pushing the existing rows 151-167 down to 152-168.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151 (shifts rows 151..167 -> 152..168)
$ws.Rows.Item(151).Insert()

# Fill in the new row 151 with the new weekly record
$ws.Range("A151").Value2 = 7
$ws.Range("B151").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value2 = "Ñuble"
$ws.Range("D151").Value2 = 45131
$ws.Range("E151").Value2 = 16
$ws.Range("F151").Value2 = "Fruta"
$ws.Range("G151").Value2 = 100108
$ws.Range("H151").Value2 = "Tropicales y subtropicales"
$ws.Range("I151").Value2 = 100108002
$ws.Range("J151").Value2 = "Mango"
$ws.Range("K151").Value2 = "Sin especificar"
$ws.Range("L151").Value2 = "Primera"
$ws.Range("M151").Value2 = 60
$ws.Range("N151").Value2 = 8000
$ws.Range("O151").Value2 = 8000
$ws.Range("P151").Value2 = 8000
$ws.Range("Q151").Value2 = "$/bandeja 4 kilos"
$ws.Range("R151").Value2 = "Brasil"
$ws.Range("S151").Value2 = 2000
$ws.Range("T151").Value2 = 4
